# Update Price (D) and Volume(1h) (E) columns for rows 2-51 per latest crypto data snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.152.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -2.59%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''1.866.97'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -2.13%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  -0.22%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''306.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.10%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.19%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''0.5163'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '''  +3.07%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -1.76%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.07173'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.61%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''20.71'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -1.11%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.8856'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -2.52%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''0.07578'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.10%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''1.839.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  -4.17%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''5.335'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.71%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = '''  -2.49%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''1.002'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.14%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''0.000008565'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.77%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''14.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -2.67%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''0.9999'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = '''27.172.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.67%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''5.042'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -2.56%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''2.093.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -5.70%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''10.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -2.10%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''6.475'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -1.93%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''150.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -2.17%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.848'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -1.70%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''18.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -2.14%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.132'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -4.29%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''112.71'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  -2.27%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''4.752'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -3.31%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''4.691'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  +0.68%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''0.08999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  +0.23%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.05159'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -1.78%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''3.096'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -3.11%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.7519'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -1.94%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''1.173'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -4.92%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''0.02038'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -1.16%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''2.543'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -0.55%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''3.031'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  +0.35%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''1.079'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -1.26%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''0.5353'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -3.88%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''6.645'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -4.59%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''114.48'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +2.59%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''8.484'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  -0.29%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''0.1483'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -1.86%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''0.4671'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -3.13%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.9993'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -0.17%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''10.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -4.34%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -3.65%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''64.99'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -3.93%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''36.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -1.56%  '
$ws.Range("E51").Style = "Normal"
